$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New combined ligand-receptor table (Rspo3-Lgr6) across all three sending/target
# cluster pairs (ECs, FAPs, sCs), per Dr Hou's advice.
$data = @(
    @("ECs",  "Rspo3", "Lgr6", "ECs", 1, 0.3333333333333333, 0.06762866666666667, 0.202886,   0.06580093600101189, 0.06580093600101189, 1, 0.3333333333333333, 0.01182833333333333, 0.035485,  0.01157399102261438, 0.01157399102261438, 0.0007999344122222224, 0.007199409710000001, 0.0007615794425553348, 0.0007615794425553347),
    @("ECs",  "Rspo3", "Lgr6", "sCs", 1, 0.3333333333333333, 0.06762866666666667, 0.202886,   0.06580093600101189, 0.06580093600101189, 3, 1,                  1.010147,             3.030441,  0.9884260089773856,  0.9884260089773856,  0.06831489474733334,  0.614834052726,       0.06503935655845655,  0.06503935655845655),
    @("FAPs", "Rspo3", "Lgr6", "ECs", 2, 0.6666666666666666, 0.3764723333333334,  1.129417,   0.3662978014030285,  0.3662978014030286,  1, 0.3333333333333333, 0.01182833333333333, 0.035485,  0.01157399102261438, 0.01157399102261438, 0.004453040249444445, 0.040077362245,       0.004239527465042037, 0.004239527465042037),
    @("FAPs", "Rspo3", "Lgr6", "sCs", 2, 0.6666666666666666, 0.3764723333333334,  1.129417,   0.3662978014030285,  0.3662978014030286,  3, 1,                  1.010147,             3.030441,  0.9884260089773856,  0.9884260089773856,  0.3802923980996666,   3.422631582897,       0.3620582739379865,   0.3620582739379866),
    @("sCs",  "Rspo3", "Lgr6", "ECs", 3, 1,                  0.5836756666666667,  1.751027,   0.5679012625959595,  0.5679012625959595,  1, 0.3333333333333333, 0.01182833333333333, 0.035485,  0.01157399102261438, 0.01157399102261438, 0.006903910343888891, 0.06213519309500001,  0.006572884115017006, 0.006572884115017005),
    @("sCs",  "Rspo3", "Lgr6", "sCs", 3, 1,                  0.5836756666666667,  1.751027,   0.5679012625959595,  0.5679012625959595,  3, 1,                  1.010147,             3.030441,  0.9884260089773856,  0.9884260089773856,  0.5895982236563333,   5.306384012906999,    0.5613283784809425,   0.5613283784809425)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
